$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master_Tracking")

# Row 2 updates (columns H-S)
$ws.Range("H2").Value = 0.0155
$ws.Range("I2").Value = 0.0083
$ws.Range("J2").Value = 0.993
$ws.Range("K2").Value = 0.0057
$ws.Range("L2").Value = 0.9931
$ws.Range("M2").Value = 0.0056
$ws.Range("N2").Value = 0.9863
$ws.Range("O2").Value = 0.0111
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.993
$ws.Range("S2").Value = 0.0057

# Row 3 updates (columns H-S)
$ws.Range("H3").Value = 0.0139
$ws.Range("I3").Value = 0.0071
$ws.Range("J3").Value = 0.994
$ws.Range("K3").Value = 0.0046
$ws.Range("L3").Value = 0.9941
$ws.Range("M3").Value = 0.0045
$ws.Range("N3").Value = 0.9882
$ws.Range("O3").Value = 0.0089
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.994
$ws.Range("S3").Value = 0.0046
